$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: numeroDocumento / usuario updated to new test account number
$ws.Range("B3").Value = "'10757881"
$ws.Range("D3").Value = "'10757881"
# Row 3: nuevoUser updated to new test username
$ws.Range("M3").Value = "'pruebasqa94"

# Row 4: usuario corrected to match numeroDocumento (was a stale username)
$ws.Range("D4").Value = "'95400152"

# Re-apply font to B4/D4 (font size stays 12, but is refreshed)
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("D4").Font.Name = "Calibri"

# Selection moved to C10
$ws.Activate()
$ws.Range("C10").Select()
